# Add a new "time_taken" column (F) with metadata timestamps to the
# "Gene therapy clinical trials" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled the same way as the other header cells (bold/centered).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row timestamps recorded when the panel data was generated.
$timestamps = @(
    "2021-10-05 13:39:55.443678",
    "2021-10-05 13:39:55.443691",
    "2021-10-05 13:39:55.443695",
    "2021-10-05 13:39:55.443698",
    "2021-10-05 13:39:55.443701",
    "2021-10-05 13:39:55.443705",
    "2021-10-05 13:39:55.443708",
    "2021-10-05 13:39:55.443711",
    "2021-10-05 13:39:55.443714",
    "2021-10-05 13:39:55.443717",
    "2021-10-05 13:39:55.443720",
    "2021-10-05 13:39:55.443723",
    "2021-10-05 13:39:55.443726",
    "2021-10-05 13:39:55.443729",
    "2021-10-05 13:39:55.443732",
    "2021-10-05 13:39:55.443735",
    "2021-10-05 13:39:55.443738",
    "2021-10-05 13:39:55.443741",
    "2021-10-05 13:39:55.443744",
    "2021-10-05 13:39:55.443747",
    "2021-10-05 13:39:55.443750"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
